# Insert a new column before A, shifting the existing "nombre"/"criterio"
# table (A:B) to (B:C). This preserves the header styling (bold, thin
# border, centered/top alignment) and the shared text values, which all
# move one column to the right automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Insert()

# Give the new index column (A2:A6) the same formatting as the header
# cells (bold / bordered / centered-top), matching style index 1 used
# elsewhere in the sheet.
$ws.Range("B1").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new zero-based index column.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# Update the "criterio" counts (now in column C) for the rows whose
# values changed.
$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 12
